# Weekly data refresh: a new observation is inserted as row 333 (above the
# existing row 333), pushing the former rows 333-372 down to 334-373.
# The sheet's used-range / dimension grows from A1:T372 to A1:T373.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 333 (shifts 333.. down by one,
# inherits formatting from the row above as Excel normally does).
$ws.Rows.Item(333).Insert()

# Populate the newly inserted row 333 with the new weekly record.
$ws.Range("A333").Value = 3
$ws.Range("B333").Value = 'Femacal de La Calera'
$ws.Range("C333").Value = 'Coquimbo'
$ws.Range("D333").Value = 45180
$ws.Range("E333").Value = 5
$ws.Range("F333").Value = 'Fruta'
$ws.Range("G333").Value = 100101
$ws.Range("H333").Value = 'Berries'
$ws.Range("I333").Value = 100101001
$ws.Range("J333").Value = 'Arándano (blue)'
$ws.Range("K333").Value = 'Sin especificar'
$ws.Range("L333").Value = 'Primera'
$ws.Range("M333").Value = 45
$ws.Range("N333").Value = 13000
$ws.Range("O333").Value = 13000
$ws.Range("P333").Value = 13000
$ws.Range("Q333").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R333").Value = 'Provincia de Quillota'
$ws.Range("S333").Value = 8667
$ws.Range("T333").Value = 1.5
